# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2-252) from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value = 45188
    }
}
